$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 7038.7393
$ws.Range("I64").Value = 5462.8
$ws.Range("K64").Value = 5462.8
$ws.Range("M64").Value = -5214.8
$ws.Range("H67").Value = 7038.7393
$ws.Range("I67").Value = 5462.8
$ws.Range("K67").Value = 5462.8
$ws.Range("M67").Value = -4604.8
$ws.Range("H107").Value = 237.33333
$ws.Range("I107").Value = 231.68182
$ws.Range("J107").Value = 299.5
$ws.Range("K107").Value = 231.68182
$ws.Range("L107").Value = 299.5
$ws.Range("M107").Value = 1688.31818
$ws.Range("N107").Value = -4139.5
$ws.Range("H127").Value = 844.75
$ws.Range("I127").Value = 844.75
$ws.Range("K127").Value = 2534.25
$ws.Range("M127").Value = 2425.75
$ws.Range("H138").Value = 6877.8823
$ws.Range("J138").Value = 7335.88
$ws.Range("L138").Value = 22007.64
$ws.Range("N138").Value = -32287.64

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20230.727
$ws.Range("I32").Value = 19835.441
$ws.Range("K32").Value = 19835.441
$ws.Range("M32").Value = -19548.441
$ws.Range("H44").Value = 34000
$ws.Range("J44").Value = 34000
$ws.Range("L44").Value = 34000
$ws.Range("N44").Value = -34976
$ws.Range("H45").Value = 3973.0588
$ws.Range("I45").Value = 3384
$ws.Range("K45").Value = 3384
$ws.Range("M45").Value = -3007
$ws.Range("H61").Value = 7806.5625
$ws.Range("I61").Value = 8349.200000000001
$ws.Range("J61").Value = 6902.1665
$ws.Range("K61").Value = 8349.200000000001
$ws.Range("L61").Value = 6902.1665
$ws.Range("M61").Value = -8137.200000000001
$ws.Range("N61").Value = -7326.1665
$ws.Range("H97").Value = 1159.2812
$ws.Range("J97").Value = 1515.3334
$ws.Range("L97").Value = 1515.3334
$ws.Range("N97").Value = -2507.3334
$ws.Range("H110").Value = 3321.2778
$ws.Range("I110").Value = 1646.6428
$ws.Range("K110").Value = 1646.6428
$ws.Range("M110").Value = 398.3571999999999
$ws.Range("H122").Value = 3035.8
$ws.Range("I122").Value = 3045
$ws.Range("K122").Value = 9135
$ws.Range("M122").Value = -6685
$ws.Range("H132").Value = 4010.5
$ws.Range("I132").Value = 4010.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12031.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9501.5
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 7806.5625
$ws.Range("I136").Value = 8349.200000000001
$ws.Range("J136").Value = 6902.1665
$ws.Range("K136").Value = 25047.6
$ws.Range("L136").Value = 20706.4995
$ws.Range("M136").Value = -22497.6
$ws.Range("N136").Value = -25806.4995

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 1495.6
$ws.Range("I54").Value = 1495.6
$ws.Range("K54").Value = 1495.6
$ws.Range("M54").Value = -1011.6
$ws.Range("H105").Value = 1991.3158
$ws.Range("I105").Value = 1824.1666
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 1824.1666
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -77.16660000000002
$ws.Range("N105").Value = -8494
$ws.Range("H107").Value = 2229.7693
$ws.Range("I107").Value = 1590.6364
$ws.Range("K107").Value = 1590.6364
$ws.Range("M107").Value = 329.3635999999999
$ws.Range("H134").Value = 5912.647
$ws.Range("I134").Value = 5565.8
$ws.Range("K134").Value = 16697.4
$ws.Range("M134").Value = -14162.4

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 15000
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H58").Value = 4188.793
$ws.Range("I58").Value = 2578.1
$ws.Range("K58").Value = 2578.1
$ws.Range("M58").Value = -2375.1
$ws.Range("H99").Value = 7397.8184
$ws.Range("I99").Value = 5931.3335
$ws.Range("J99").Value = 7947.75
$ws.Range("K99").Value = 5931.3335
$ws.Range("L99").Value = 7947.75
$ws.Range("M99").Value = -4433.3335
$ws.Range("N99").Value = -10943.75
$ws.Range("H126").Value = 7397.8184
$ws.Range("I126").Value = 5931.3335
$ws.Range("J126").Value = 7947.75
$ws.Range("K126").Value = 17794.0005
$ws.Range("L126").Value = 23843.25
$ws.Range("M126").Value = -15324.0005
$ws.Range("N126").Value = -28783.25
$ws.Range("H132").Value = 2443.6667
$ws.Range("I132").Value = 1602.8518
$ws.Range("K132").Value = 4808.555399999999
$ws.Range("M132").Value = -2278.555399999999
$ws.Range("H134").Value = 3234.6135
$ws.Range("I134").Value = 1566.8334
$ws.Range("K134").Value = 4700.5002
$ws.Range("M134").Value = -2165.5002
$ws.Range("H136").Value = 4188.793
$ws.Range("I136").Value = 2578.1
$ws.Range("K136").Value = 7734.299999999999
$ws.Range("M136").Value = -5184.299999999999

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 645.2727
$ws.Range("J46").Value = 699.75
$ws.Range("L46").Value = 2099.25
$ws.Range("N46").Value = -2281.25
$ws.Range("H122").Value = 16006.363
$ws.Range("I122").Value = 6412.25
$ws.Range("J122").Value = 18138.389
$ws.Range("K122").Value = 57710.25
$ws.Range("L122").Value = 163245.501
$ws.Range("M122").Value = -55260.25
$ws.Range("N122").Value = -168145.501
$ws.Range("H131").Value = 3686.889
$ws.Range("I131").Value = 2043.5
$ws.Range("J131").Value = 4156.4287
$ws.Range("K131").Value = 6130.5
$ws.Range("L131").Value = 12469.2861
$ws.Range("M131").Value = -1090.5
$ws.Range("N131").Value = -22549.2861
$ws.Range("H132").Value = 1942.5714
$ws.Range("I132").Value = 1099
$ws.Range("J132").Value = 2280
$ws.Range("K132").Value = 9891
$ws.Range("L132").Value = 20520
$ws.Range("M132").Value = -7361
$ws.Range("N132").Value = -25580
$ws.Range("H140").Value = 387642.53
$ws.Range("I140").Value = 1910.8462
$ws.Range("K140").Value = 5732.5386
$ws.Range("M140").Value = -552.5385999999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7669.7646
$ws.Range("I102").Value = 5438.8
$ws.Range("K102").Value = 5438.8
$ws.Range("M102").Value = -3816.8
$ws.Range("H122").Value = 5206.524
$ws.Range("I122").Value = 5164.5713
$ws.Range("K122").Value = 15493.7139
$ws.Range("M122").Value = -13043.7139
$ws.Range("H132").Value = 3199.1538
$ws.Range("I132").Value = 2539.1428
$ws.Range("J132").Value = 5971.2
$ws.Range("K132").Value = 7617.428400000001
$ws.Range("L132").Value = 17913.6
$ws.Range("M132").Value = -5087.428400000001
$ws.Range("N132").Value = -22973.6

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 16808.63
$ws.Range("I132").Value = 16221.3125
$ws.Range("J132").Value = 19941
$ws.Range("K132").Value = 48663.9375
$ws.Range("L132").Value = 59823
$ws.Range("M132").Value = -46133.9375
$ws.Range("N132").Value = -64883

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H107").Value = 1557.878
$ws.Range("I107").Value = 1220.04
$ws.Range("K107").Value = 3660.12
$ws.Range("M107").Value = -1740.12
$ws.Range("H113").Value = 4184.467
$ws.Range("I113").Value = 1751.4615
$ws.Range("K113").Value = 5254.3845
$ws.Range("M113").Value = -3084.3845
$ws.Range("H132").Value = 3882.8645
$ws.Range("I132").Value = 3208.0889
$ws.Range("K132").Value = 9624.2667
$ws.Range("M132").Value = -7094.2667
$ws.Range("H136").Value = 4435.533
$ws.Range("I136").Value = 2849.842
$ws.Range("K136").Value = 8549.526
$ws.Range("M136").Value = -5999.526
